$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'29.957.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.893.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.17%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - XRP
$ws.Range("D5").Value = "'0.7757"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'243.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3131"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.34%  "

# Row 9 - Solana
$ws.Range("E9").Value = "  +2.03%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.07243"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.40%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.08719"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.98%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "'2.075.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.18%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "'0.7748"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.421"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.42%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'94.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.42%  "

# Row 16 - Uniswap (was WrappedBTC)
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "'6.190"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.73%  "

# Row 17 - WrappedBTC (was Uniswap)
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'30.184.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.17%  "

# Row 18 - WrappedliquidstakedEther2.0
$ws.Range("D18").Value = "'2.360.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.18%  "

# Row 19 - Avalanche
$ws.Range("E19").Value = "  -0.06%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'245.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.09%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "'0.000007860"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.34%  "

# Row 22 - Dai
$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "'8.114"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  -0.10%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "'0.1644"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.93%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'9.485"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.00%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'163.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +0.54%  "

# Row 29 - LidoDAOToken
$ws.Range("E29").Value = "  +0.73%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -0.20%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.44%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.518"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.19%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.142"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.05475"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "'1.244"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.22%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'0.7538"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.06%  "

# Row 37 - Frax
$ws.Range("D37").Value = "'1.008"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.60%  "

# Row 38 - HuobiToken
$ws.Range("D38").Value = "'2.701"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.93%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +2.97%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "'2.788"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.36%  "

# Row 41 - TheSandbox
$ws.Range("E41").Value = "  +2.33%  "

# Row 42 - Maker
$ws.Range("D42").Value = "'1.111.19"
$ws.Range("D42").Style = "Normal"

# Row 43 - Aave (was FraxShare)
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'73.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.37%  "

# Row 44 - FraxShare (was Aave)
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.118"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.89%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "'2.243.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.22%  "

# Row 46 - TrustWalletToken
$ws.Range("D46").Value = "'0.8497"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47 - Quant
$ws.Range("D47").Value = "'104.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "

# Row 49 - RenderToken
$ws.Range("D49").Value = "'1.879"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "

# Row 50 - Aptos
$ws.Range("D50").Value = "'7.624"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.37%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "'9.865"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.34%  "

